$wb = $excel.ActiveWorkbook

# --- "axes" sheet: add arrow-label and percent-label columns ------------
$axes = $wb.Worksheets.Item("axes")

# Row 1 (headers): keep A/B/C, insert the three "*_arrow" headers before
# the existing "Title" column (which slides from D1 to G1).
$axes.Range("D1").Value = "A_arrow"
$axes.Range("E1").Value = "B_arrow"
$axes.Range("F1").Value = "C_arrow"
$axes.Range("G1").Value = "Title"

# Row 2 (values): keep Q/F/L, insert the three "* (%)" values before the
# existing "QFL Diagram" column (which slides from D2 to G2).
$axes.Range("D2").Value = "Q (%)"
$axes.Range("E2").Value = "F (%)"
$axes.Range("F2").Value = "L (%)"
$axes.Range("G2").Value = "QFL Diagram"

# Widen the three new columns (target stored width ~14.11 chars; the
# engine snaps ColumnWidth assignments to its pixel grid, so 13.333 lands
# on the closest reachable stored width).
$axes.Range("D1:F1").ColumnWidth = 13.333

# Move the selection to F8 and make "axes" the active/visible tab (this
# also clears tabSelected on whichever sheet had it before).
$axes.Range("F8").Select()
$axes.Activate()
